$d = $word.ActiveDocument

# 1. Merge the title runs "INFORME DE LA PRUEBA" + " DE REGRESIÓN"
#    into a single run "INFORME DE LA PRUEBA DE REGRESIÓN".
$d.Content.Find.Execute("INFORME DE LA PRUEBA DE REGRESIÓN", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "INFORME DE LA PRUEBA DE REGRESIÓN", 2)

# 2. Change "Hoja 1/1" to "1 Hoja"
$d.Content.Find.Execute("Hoja 1/1", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1 Hoja", 2)
